$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (H1) onto the two
# new header cells so they pick up the same style index (bold/border/
# centered header look) instead of Excel creating a brand new style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data row values (plain, unstyled numeric cells like the rest of row 2)
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
